# Saldo.xlsx update
#
# Net effect of the target diff on the "Export" sheet:
#   - 5 account rows get an updated Saldo (column C) value
#   - 18 account rows are removed entirely
#   - the data block (still header row 1, then a blank row + a filter-notes
#     row at the bottom) stays sorted descending by Saldo, so after the
#     value edits / deletions we re-sort the data range to restore that
#     invariant (this is what naturally re-positions rows such as WILSON's).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlWhole lookAt constant used with Range.Find
$xlWhole = 1

function Find-AccountRow([string]$account) {
    $lastRow = $ws.UsedRange.Rows.Count
    $colA = $ws.Range("A1:A" + $lastRow)
    $hit = $colA.Find($account, [Type]::Missing, [Type]::Missing, $xlWhole)
    if ($null -eq $hit) {
        return -1
    }
    return $hit.Row
}

# 1) Update Saldo values for accounts whose balance changed.
$changedValues = @{
    "004641487" = 308711.41   # LAILA        306007.86 -> 308711.41
    "004467884" = 3120.71     # ANA          8120.71   -> 3120.71
    "005320069" = 2132.44     # RICARDO      7981.25   -> 2132.44
    "004397124" = 723.32      # MURYLO       7023.32   -> 723.32
    "004884046" = 1126.4      # WILSON       499.15    -> 1126.4
}

foreach ($account in $changedValues.Keys) {
    $row = Find-AccountRow $account
    if ($row -gt 0) {
        $ws.Cells.Item($row, 3).Value = $changedValues[$account]
    }
}

# 2) Remove accounts that no longer appear in the export.
$removedAccounts = @(
    "004472404",  # DILSON     13685.24
    "005122672",  # LUCAS      8500
    "004450724",  # ASSAKO     8474.84
    "005366671",  # TATIANA    6633.86
    "004514241",  # ANDRE      6237.2
    "004331477",  # SUZY       5518.03
    "001294033",  # VIVIANE    5140.61
    "004267044",  # PATRICIA   4413.24
    "004384167",  # DOUGLAS    4090.61
    "004384258",  # PAULA      3820.45
    "004482090",  # CEZAR      2843.41
    "004265173",  # JULIA      2431.59
    "004374943",  # LEONARDO   2203.46
    "005046919",  # MARIANA    2004.24
    "004643737",  # LARA       1868.66
    "004855596",  # MARIANA    1792.49
    "004240400",  # ADRIANA    1082.61
    "001000882"   # AYRTON     810.12
)

foreach ($account in $removedAccounts) {
    $row = Find-AccountRow $account
    if ($row -gt 0) {
        $ws.Rows($row).Delete()
    }
}

# 3) Re-sort the data (rows 2..lastDataRow) descending by Saldo (column C)
#    to restore the sheet's sort invariant after the edits above. The sheet
#    always ends with one blank row followed by a "Filtros aplicados" notes
#    row, so the data block ends two rows above the last used row.
$lastUsedRow = $ws.UsedRange.Rows.Count
$lastDataRow = $lastUsedRow - 2

$dataRange = $ws.Range("A2:C" + $lastDataRow)
$keyRange = $ws.Range("C2:C" + $lastDataRow)
$dataRange.Sort($keyRange, 2)
